# Generate Report for Handback
# Adds a new row (row 4) describing the handback status of the file
# "e53b230d-62dc-47d3-bd27-7ad03d8c94e3.md" to the Overview, zh-cn and
# de-de worksheets, mirroring the existing rows for the other two files.

$wb = $excel.ActiveWorkbook

$guid = "e53b230d-62dc-47d3-bd27-7ad03d8c94e3"
$mdName = "$guid.md"
$inSync = "Handed back: in sync with en-US"

# Helper colour matching the workbook's custom "HyperLink" cell style
# (font color FF6495ED, underlined, Calibri 11) used throughout the sheet.
$hyperlinkColor = 15570276

function Set-HyperlinkCell($ws, $cellRef, $url, $text) {
    $rng = $ws.Range($cellRef)
    $ws.Hyperlinks.Add($rng, $url, "", "", $text) | Out-Null
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = $true
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-HyperlinkCell $wsOverview "A4" `
    "https://github.com/OpenLocalizationTest/oltest/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/e2e/$mdName" `
    $mdName
$wsOverview.Range("B4").Value = $inSync
$wsOverview.Range("C4").Value = $inSync

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$guid.6af7272b3b71e900f86fa71e7987c7bd8ab31517.zh-cn.xlf"

Set-HyperlinkCell $wsZh "A4" `
    "https://github.com/OpenLocalizationTest/oltest/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/e2e/$mdName" `
    $mdName
Set-HyperlinkCell $wsZh "B4" `
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/e2e/$mdName" `
    ".md"
$wsZh.Range("C4").Value = $inSync
Set-HyperlinkCell $wsZh "D4" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf" `
    $zhXlf
$wsZh.Range("E4").Value = "2016-03-17 02:49:36"
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-HyperlinkCell $wsZh "F4" `
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/e2e/$mdName" `
    $mdName
Set-HyperlinkCell $wsZh "G4" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf" `
    $zhXlf
$wsZh.Range("H4").Value = "2016-03-17 02:50:17"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf = "$guid.6af7272b3b71e900f86fa71e7987c7bd8ab31517.de-de.xlf"

Set-HyperlinkCell $wsDe "A4" `
    "https://github.com/OpenLocalizationTest/oltest/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/e2e/$mdName" `
    $mdName
Set-HyperlinkCell $wsDe "B4" `
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/e2e/$mdName" `
    ".md"
$wsDe.Range("C4").Value = $inSync
Set-HyperlinkCell $wsDe "D4" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf" `
    $deXlf
$wsDe.Range("E4").Value = "2016-03-17 02:49:44"
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-HyperlinkCell $wsDe "F4" `
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/e2e/$mdName" `
    $mdName
Set-HyperlinkCell $wsDe "G4" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e53b230d62dc47d3bd277ad03d8c94e3000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf" `
    $deXlf
$wsDe.Range("H4").Value = "2016-03-17 02:50:31"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = "Include"
